$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data from columns B..K on row 1 (only A1 survives)
$ws.Range("B1:K1").ClearContents()

# Row 1 / A1 now holds the new product id
$ws.Range("A1").Value = "K9F1G08U0B-PIB0"
